$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.780.68"
$ws.Range("E2").Value = "  +5.11%  "

$ws.Range("D3").Value = "2.600.34"
$ws.Range("E3").Value = "  +2.73%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "522.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.08%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.57%  "

$ws.Range("E7").Value = "  -0.30%  "

$ws.Range("E8").Value = "  +1.99%  "

$ws.Range("D9").Value = "2.623.89"
$ws.Range("E9").Value = "  +3.49%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.56"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.45%  "

$ws.Range("E11").Value = "  +2.32%  "

$ws.Range("E12").Value = "  +2.69%  "

$ws.Range("E13").Value = "  +2.00%  "

$ws.Range("D14").Value = "3.062.93"
$ws.Range("E14").Value = "  +2.82%  "

$ws.Range("D15").Value = "59.609.04"
$ws.Range("E15").Value = "  +4.75%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.45"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.50%  "

$ws.Range("D17").Value = "2.624.55"
$ws.Range("E17").Value = "  +5.27%  "

$ws.Range("E18").Value = "  +0.61%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "339.23"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.33%  "

$ws.Range("E20").Value = "  +1.60%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.65%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.54"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.99%  "

$ws.Range("E23").Value = "  -0.29%  "

$ws.Range("E24").Value = "  +3.72%  "

$ws.Range("E25").Value = "  +2.19%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.405"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.37%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.996"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.65%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.08"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.52%  "

$ws.Range("E29").Value = "  -0.17%  "

$ws.Range("D30").Value = "0.0₃0727"
$ws.Range("E30").Value = "  -2.75%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.96"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.78%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.82"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.09%  "

$ws.Range("E33").Value = "  +1.39%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "149.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.39%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.02"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.72%  "

$ws.Range("E36").Value = "  +0.52%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "36.38"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.12%  "

$ws.Range("E38").Value = "  +4.39%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.833"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.96%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.823"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.95%  "

$ws.Range("E41").Value = "  +2.63%  "

$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "277.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.64%  "

$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.998"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.30%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.76"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.61%  "

$ws.Range("E45").Value = "  +3.55%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0953"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.28%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "18.67"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.67%  "

$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0521"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.38%  "

$ws.Range("D49").Value = "1.986.72"
$ws.Range("E49").Value = "  +1.18%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.62"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.60%  "

$ws.Range("E51").Value = "  +0.23%  "
